$wb = $excel.ActiveWorkbook

# --- Rename sheets (pseudotime early + late cluster relabeling) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Early_C1,2down"
$ws2.Name = "Early_C4,5down"

# --- Sheet "Early_C1,2down": hide helper columns D:F and H ---
$ws1.Columns("D:F").ColumnWidth = -0.8333333333333334
$ws1.Columns("D:F").Hidden = $true
$ws1.Columns("H").ColumnWidth = -0.8333333333333334
$ws1.Columns("H").Hidden = $true

# Move the yellow "picked" highlight from C9 to C2
$ws1.Range("C2").Interior.Color = 65535
$ws1.Range("C9").Interior.Pattern = -4142

# --- Update the remembered selection on each sheet ---
# (select sheet2 / sheet3 first so the workbook ends up back on sheet1,
#  matching the original tabSelected="1" on Early_C1,2down)
$ws2.Range("I5").Select() | Out-Null
$ws3.Range("C2").Select() | Out-Null
$ws1.Range("I47").Select() | Out-Null
